# Commit: "fix single quoted sheet names"
# Rename the sheet "This is a sheet" to "This is a - sheet" (adds a
# hyphen so the name exercises the single-quoted sheet-name syntax in
# formula references, e.g. 'This is a - sheet'!$B$2). Excel will
# automatically update any formulas that reference the sheet, and the
# renamed sheet becomes the active tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("This is a sheet")
$ws1.Name = "This is a - sheet"

# The renamed sheet becomes the active/selected tab.
$ws1.Activate()
